# Re-export of the KiBot BoM after a schematic change:
#  - the 0 ohm resistor R2 group is gone
#  - the connector reference changed from J2 (edge) to J1 (middle),
#    which used to be the DNF'd connector and is now the fitted one
#  - component counters are updated to match

$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")
$dnf = $wb.Worksheets.Item("DNF")

# --- BoM sheet -------------------------------------------------------

# Summary counters (rows 2-6, column F)
$bom.Range("F2").Value = 5
$bom.Range("F3").Value = "27 (27 SMD/ 0 THT)"
$bom.Range("F4").Value = "26 (26 SMD/ 0 THT)"
$bom.Range("F6").Value = 26

# Row 11: connector group J2/edge -> J1/middle
$bom.Range("D11").Value = "J1"
$bom.Range("E11").Value = "middle"

# Row 12: resistor group R2/0 -> R1/330 (with updated datasheet/supplier links).
# E12's new text ("330") looks numeric, so a plain .Value assignment would
# silently turn the cell into a number; pull the text over (as text, via
# xlPasteValues) from row 13, which already holds "330" as a text value,
# so the cell keeps its original text type and style.
$bom.Range("D12").Value = "R1"
$bom.Range("E13").Copy() | Out-Null
$bom.Range("E12").PasteSpecial(-4163) | Out-Null
$bom.Range("I12").Value = "https://api.pim.na.industrial.panasonic.com/file_stream/main/fileversion/1242"
$bom.Range("J12").Value = "https://www.digikey.ch/en/products/detail/panasonic-electronic-components/ERJ-1GNJ331C/2035775"

# Row 13 (old R1/330 group, now redundant) is removed entirely
$bom.Rows.Item(13).Delete()

# --- DNF sheet ---------------------------------------------------------

# Summary counters (rows 2-6, column F) mirror the BoM sheet
$dnf.Range("F2").Value = 5
$dnf.Range("F3").Value = "27 (27 SMD/ 0 THT)"
$dnf.Range("F4").Value = "26 (26 SMD/ 0 THT)"
$dnf.Range("F6").Value = 26

# Row 10 (connector J1/middle) is no longer DNF, remove it
$dnf.Rows.Item(10).Delete()

# Column widths re-fit to the now-shorter remaining content
$dnf.Columns.Item(2).ColumnWidth = 44.7109375
$dnf.Columns.Item(3).ColumnWidth = 15.7109375
$dnf.Columns.Item(6).ColumnWidth = 19.7109375
$dnf.Columns.Item(9).ColumnWidth = 19.7109375
